$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows: swap the person names for the new set, and
# refresh the city names in column C.
$ws.Range("A1").Value = "Jokowi"
$ws.Range("B1").Value = 20202
$ws.Range("C1").Value = "Barabai"

$ws.Range("A2").Value = "Prabowo"
$ws.Range("B2").Value = 76767
$ws.Range("C2").Value = "Kandangan"

# Add two brand-new rows of data.
$ws.Range("A3").Value = "Gibran"
$ws.Range("B3").Value = 99002
$ws.Range("C3").Value = "Banjarmasin"

$ws.Range("A4").Value = "Ganjar"
$ws.Range("B4").Value = 85621
$ws.Range("C4").Value = "Tabalong"

# Move the active selection below the new data, matching the saved file.
$ws.Range("A5").Select() | Out-Null
